$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F5").Value = "커트시를 하는 듯 우아한 움직임으로 사용하는 고양이 일족의 무기."
$ws.Range("F6").Value = "토끼 농부들이 사용하는 낫. 낮에는 밭을 매고 밤에는 수렵하는... 그들은 '노동자'다."
$ws.Range("F19").Value = "어린 토끼들이 자주 사용하는 호미. 낮에는 밭을 매고 밤에는 수렵하는… 이런 말 안듣는 아이들…"
